# Update the "Data" sheet (invalid drugname test case) with new sample values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Row 2: Transfer out Imprest -> Transfer out Patient, Acetec 5 mg tablet -> hmhhm
$ws.Range("A2").Value = "Transfer out Patient"
$ws.Range("C2").Value = "hmhhm"

# Row 3: Transfer out Imprest -> Transfer out Patient, Endone 5 mg tablet -> mhmhm
$ws.Range("A3").Value = "Transfer out Patient"
$ws.Range("C3").Value = "mhmhm"

# Leave the last active cell on A8, matching the author's final selection.
$ws.Range("A8").Select()
